$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shifted weekly records (rows 80-119): Date (D), and the few Volume (J) / Unit (N) values that move with them ---
$ws.Cells.Item(80, 4).Value = 44553
$ws.Cells.Item(81, 4).Value = 44553
$ws.Cells.Item(82, 4).Value = 44285
$ws.Cells.Item(83, 4).Value = 44285
$ws.Cells.Item(84, 4).Value = 44160
$ws.Cells.Item(85, 4).Value = 44160
$ws.Cells.Item(86, 4).Value = 44231
$ws.Cells.Item(87, 4).Value = 44231
$ws.Cells.Item(88, 4).Value = 44490
$ws.Cells.Item(89, 4).Value = 44490
$ws.Cells.Item(90, 4).Value = 44341
$ws.Cells.Item(91, 4).Value = 44341
$ws.Cells.Item(92, 4).Value = 44391
$ws.Cells.Item(93, 4).Value = 44391
$ws.Cells.Item(94, 4).Value = 44386
$ws.Cells.Item(94, 10).Value = 200
$ws.Cells.Item(95, 4).Value = 44386
$ws.Cells.Item(95, 10).Value = 100
$ws.Cells.Item(96, 4).Value = 44278
$ws.Cells.Item(96, 10).Value = 300
$ws.Cells.Item(97, 4).Value = 44278
$ws.Cells.Item(97, 10).Value = 150
$ws.Cells.Item(98, 4).Value = 44308
$ws.Cells.Item(99, 4).Value = 44308
$ws.Cells.Item(100, 4).Value = 44187
$ws.Cells.Item(101, 4).Value = 44187
$ws.Cells.Item(102, 4).Value = 44350
$ws.Cells.Item(103, 4).Value = 44350
$ws.Cells.Item(104, 4).Value = 44405
$ws.Cells.Item(105, 4).Value = 44405
$ws.Cells.Item(106, 4).Value = 44224
$ws.Cells.Item(107, 4).Value = 44224
$ws.Cells.Item(108, 4).Value = 44398
$ws.Cells.Item(109, 4).Value = 44398
$ws.Cells.Item(110, 4).Value = 44239
$ws.Cells.Item(110, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(111, 4).Value = 44239
$ws.Cells.Item(111, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(112, 4).Value = 44344
$ws.Cells.Item(112, 14).Value = "`$/docena de 1 kilo"
$ws.Cells.Item(113, 4).Value = 44344
$ws.Cells.Item(113, 14).Value = "`$/docena de 1 kilo"
$ws.Cells.Item(114, 4).Value = 44365
$ws.Cells.Item(115, 4).Value = 44365
$ws.Cells.Item(116, 4).Value = 44194
$ws.Cells.Item(117, 4).Value = 44194
$ws.Cells.Item(118, 4).Value = 44313
$ws.Cells.Item(119, 4).Value = 44313

# --- Append two new rows (120, 121) duplicating the oldest pair of records that got pushed off ---
$ws.Cells.Item(120, 1).Value = 11
$ws.Cells.Item(120, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value = "Bíobío"
$ws.Cells.Item(120, 4).Value = 44272
$ws.Cells.Item(120, 5).Value = 8
$ws.Cells.Item(120, 6).Value = 100112044
$ws.Cells.Item(120, 7).Value = "Perejil"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 200
$ws.Cells.Item(120, 11).Value = 600
$ws.Cells.Item(120, 12).Value = 700
$ws.Cells.Item(120, 13).Value = 650
$ws.Cells.Item(120, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(120, 15).Value = "Región de Ñuble"
$ws.Cells.Item(120, 16).Value = 650
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
$ws.Cells.Item(120, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(121, 1).Value = 11
$ws.Cells.Item(121, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(121, 3).Value = "Bíobío"
$ws.Cells.Item(121, 4).Value = 44272
$ws.Cells.Item(121, 5).Value = 8
$ws.Cells.Item(121, 6).Value = 100112044
$ws.Cells.Item(121, 7).Value = "Perejil"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Segunda"
$ws.Cells.Item(121, 10).Value = 100
$ws.Cells.Item(121, 11).Value = 500
$ws.Cells.Item(121, 12).Value = 500
$ws.Cells.Item(121, 13).Value = 500
$ws.Cells.Item(121, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(121, 15).Value = "Región de Ñuble"
$ws.Cells.Item(121, 16).Value = 500
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
